$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "LAST SCRAPE DATE" column (F) everywhere it currently reads
#    2019-03-07 so it reads 2019-03-12 instead. The value must remain a
#    text string (not get converted into a numeric date), so force a text
#    number format before assigning the new value.
$oldDate = "2019-03-07"
$newDate = "2019-03-12"
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq $oldDate) {
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
    }
}

# 2) Swap the "Snow Me The Money" / "White Hot 7's" rows' data (game name,
#    game number, top prizes remaining) between rows 6 and 7.
$c6 = $ws.Range("C6").Value2
$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2
$c7 = $ws.Range("C7").Value2
$d7 = $ws.Range("D7").Value2
$e7 = $ws.Range("E7").Value2

$ws.Range("C6").Value = $c7
$ws.Range("D6").Value = $d7
$ws.Range("E6").Value = $e7
$ws.Range("C7").Value = $c6
$ws.Range("D7").Value = $d6
$ws.Range("E7").Value = $e6

# 3) Update the "TOP PRIZES REMAINING" count for the Tripler Cashword row
#    (row 27) from 2 down to 1.
$ws.Range("E27").Value = 1
